# "Update countries & provincias Spain"
# Refresh the COVID-19 country stats table (sheet "Pais") with the newer
# snapshot, update the "last updated" timestamp, and re-rank
# Vietnam/Tanzania (Vietnam's totals overtook Tanzania's, so the two swap
# table positions).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: India -------------------------------------------------------
$ws.Range("B6").Value = 1649323
$ws.Range("C6").Value = 9973
$ws.Range("D6").Value = 1063296
$ws.Range("E6").Value = 550140
$ws.Range("G6").Value = 101
$ws.Range("H6").Value = 35887

# --- Row 26: Catar -------------------------------------------------------
$ws.Range("B26").Value = 110695
$ws.Range("C26").Value = 235
$ws.Range("D26").Value = 107377
$ws.Range("E26").Value = 3144
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 174

# --- Row 40: Bielorrusia ---------------------------------------------------
$ws.Range("B40").Value = 67808
$ws.Range("C40").Value = 143
$ws.Range("D40").Value = 62444
$ws.Range("E40").Value = 4805
$ws.Range("G40").Value = 6
$ws.Range("H40").Value = 559

# --- Row 41: Kuwait ------------------------------------------------------
$ws.Range("B41").Value = 66957
$ws.Range("C41").Value = 428
$ws.Range("D41").Value = 57932
$ws.Range("E41").Value = 8578
$ws.Range("G41").Value = 2
$ws.Range("H41").Value = 447

# --- Row 52: Barein --------------------------------------------------------
$ws.Range("E52").Value = 3251
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 147

# --- Row 68: Nepal ---------------------------------------------------------
$ws.Range("B68").Value = 19771
$ws.Range("C68").Value = 224
$ws.Range("D68").Value = 14399
$ws.Range("E68").Value = 5316
$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 56

# --- Row 83: Madagascar ------------------------------------------------------
$ws.Range("B83").Value = 10868
$ws.Range("C83").Value = 120
$ws.Range("D83").Value = 7807
$ws.Range("E83").Value = 2955
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 106

# --- Row 87: Consejo Danes para los Refugiados ------------------------------
$ws.Range("B87").Value = 9070
$ws.Range("C87").Value = 60
$ws.Range("D87").Value = 6796
$ws.Range("E87").Value = 2059

# --- Row 142: Liberia --------------------------------------------------------
$ws.Range("B142").Value = 1186
$ws.Range("C142").Value = 5
$ws.Range("D142").Value = 670
$ws.Range("G142").Value = 2
$ws.Range("H142").Value = 75

# --- Row 150: Principado de Andorra ------------------------------------------
$ws.Range("B150").Value = 925
$ws.Range("C150").Value = 3
$ws.Range("D150").Value = 807
$ws.Range("E150").Value = 66

# --- Row 154: Malta -----------------------------------------------------------
$ws.Range("B154").Value = 824
$ws.Range("C154").Value = 10
$ws.Range("E154").Value = 150

# --- Rows 161/162: Tanzania & Vietnam swap ranking ---------------------------
# Vietnam's case count (545) overtook Tanzania's (509), so Vietnam moves up
# to row 161 (bringing its refreshed stats) and Tanzania drops to row 162
# (keeping its previous, unchanged stats).
$ws.Range("A161").Value = "Vietnam"
$ws.Range("B161").Value = 545
$ws.Range("C161").Value = 36
$ws.Range("D161").Value = 373
$ws.Range("E161").Value = 171
$ws.Range("G161").Value = 1
$ws.Range("H161").Value = 1

$ws.Range("A162").Value = "Tanzania"
$ws.Range("B162").Value = 509
$ws.Range("C162").Value = 0
$ws.Range("D162").Value = 183
$ws.Range("E162").Value = 305
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 21

# --- Timestamp footer ---------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 31 de Julio de 2020 a las 13:52"
